# Update the cryptocurrency price/volume table (Mon Jun  5 21:54:13 UTC 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.609.00"
$ws.Cells.Item(2, 5).Value = "  -5.89%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.804.04"
$ws.Cells.Item(3, 5).Value = "  -5.21%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "273.62"
$ws.Cells.Item(5, 5).Value = "  -10.60%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.01%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5028"
$ws.Cells.Item(7, 5).Value = "  -7.05%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.3505"
$ws.Cells.Item(8, 5).Value = "  -7.90%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "OKB"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(9, 4).Value = "43.71"
$ws.Cells.Item(9, 5).Value = "  -4.94%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value = "0.06602"
$ws.Cells.Item(10, 5).Value = "  -9.44%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "Solana"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(11, 4).Value = "19.89"
$ws.Cells.Item(11, 5).Value = "  -9.82%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Polygon"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(12, 4).Value = "0.8323"
$ws.Cells.Item(12, 5).Value = "  -7.70%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "TRON"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(13, 4).Value = "0.07776"
$ws.Cells.Item(13, 5).Value = "  -4.96%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.788.52"
$ws.Cells.Item(14, 5).Value = "  +39.55%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "Polkadot"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(15, 4).Value = "5.045"
$ws.Cells.Item(15, 5).Value = "  -5.58%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).Value = "87.30"
$ws.Cells.Item(16, 5).Value = "  -8.51%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "BinanceUSD"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(17, 4).Value = "0.9995"
$ws.Cells.Item(17, 5).Value = "  -0.11%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "Avalanche"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(18, 4).Value = "13.86"
$ws.Cells.Item(18, 5).Value = "  -6.31%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "Dai"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(19, 4).Value = "1.001"
$ws.Cells.Item(19, 5).Value = "  +0.06%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "ShibaInu"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(20, 4).Value = "0.000007920"
$ws.Cells.Item(20, 5).Value = "  -8.17%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "WrappedBTC"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(21, 4).Value = "25.659.49"
$ws.Cells.Item(21, 5).Value = "  -5.81%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(22, 4).Value = "4.697"
$ws.Cells.Item(22, 5).Value = "  -6.92%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "Cosmos"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(23, 4).Value = "9.981"
$ws.Cells.Item(23, 5).Value = "  -7.57%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "Chainlink"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(24, 4).Value = "6.042"
$ws.Cells.Item(24, 5).Value = "  -7.17%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).Value = "141.80"
$ws.Cells.Item(25, 5).Value = "  -4.41%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "LidoDAOToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(26, 4).Value = "2.108"
$ws.Cells.Item(26, 5).Value = "  -8.36%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "Toncoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(27, 4).Value = "1.651"
$ws.Cells.Item(27, 5).Value = "  -5.71%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).Value = "16.89"
$ws.Cells.Item(28, 5).Value = "  -8.00%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "BitcoinCash"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(29, 4).Value = "108.32"
$ws.Cells.Item(29, 5).Value = "  -7.20%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).Value = "4.304"
$ws.Cells.Item(30, 5).Value = "  -11.22%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).Value = "4.179"
$ws.Cells.Item(31, 5).Value = "  -10.17%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).Value = "0.08770"
$ws.Cells.Item(32, 5).Value = "  -4.69%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).Value = "0.04781"
$ws.Cells.Item(33, 5).Value = "  -5.53%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "1.129"
$ws.Cells.Item(34, 5).Value = "  -7.37%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.7182"
$ws.Cells.Item(35, 5).Value = "  -12.56%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "HuobiToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(36, 4).Value = "2.870"
$ws.Cells.Item(36, 5).Value = "  -4.70%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Frax"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(37, 4).Value = "0.9997"
$ws.Cells.Item(37, 5).Value = "  +0.03%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "MXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(38, 4).Value = "3.021"
$ws.Cells.Item(38, 5).Value = "  -8.86%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.01852"
$ws.Cells.Item(39, 5).Value = "  -7.40%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "TheSandbox"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(40, 4).Value = "0.5136"
$ws.Cells.Item(40, 5).Value = "  -14.59%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).Value = "2.269"
$ws.Cells.Item(41, 5).Value = "  -15.88%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).Value = "0.9410"
$ws.Cells.Item(42, 5).Value = "  -12.43%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).Value = "113.41"
$ws.Cells.Item(43, 5).Value = "  -2.09%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value = "6.126"
$ws.Cells.Item(44, 5).Value = "  -7.89%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(45, 4).Value = "7.986"
$ws.Cells.Item(45, 5).Value = "  -13.63%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "PaxDollar"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(46, 4).Value = "1.001"
$ws.Cells.Item(46, 5).Value = "  +0.05%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "0.4542"
$ws.Cells.Item(47, 5).Value = "  -11.99%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).Value = "0.1372"
$ws.Cells.Item(48, 5).Value = "  -10.27%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "9.275"
$ws.Cells.Item(49, 5).Value = "  -8.74%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).Value = "35.99"
$ws.Cells.Item(50, 5).Value = "  -5.51%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value = "1.478"
$ws.Cells.Item(51, 5).Value = "  -9.61%  "

